$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells remain text-formatted so values keep exact string representation
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.83"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.328"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06204"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.590"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.673"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.390"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8293"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01364"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1607"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08214"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03157"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09284"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001723"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04838"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006333"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005389"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001090"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.756"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.330"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3349"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1213"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002683"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04651"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006898"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1153"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003601"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01227"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006247"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7001"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1631"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01240"
